$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H45").Value = 703.2
$ws.Range("J45").Value = 703.2
$ws.Range("L45").Value = 2109.6
$ws.Range("N45").Value = -2493.6
$ws.Range("H98").Value = 4689.115
$ws.Range("I98").Value = 5053.6665
$ws.Range("J98").Value = 314.5
$ws.Range("K98").Value = 5053.6665
$ws.Range("L98").Value = 314.5
$ws.Range("M98").Value = -3555.6665
$ws.Range("N98").Value = -3310.5
$ws.Range("H107").Value = 1979.5161
$ws.Range("I107").Value = 1106.4615
$ws.Range("J107").Value = 6519.4
$ws.Range("K107").Value = 1106.4615
$ws.Range("L107").Value = 6519.4
$ws.Range("M107").Value = 813.5385000000001
$ws.Range("N107").Value = -10359.4
$ws.Range("H112").Value = 1957.7778
$ws.Range("J112").Value = 2203.4
$ws.Range("L112").Value = 6610.200000000001
$ws.Range("N112").Value = -8826.200000000001
$ws.Range("H116").Value = 2589.0908
$ws.Range("I116").Value = 2354.4285
$ws.Range("K116").Value = 2354.4285
$ws.Range("M116").Value = 1087.5715
$ws.Range("H122").Value = 4689.115
$ws.Range("I122").Value = 5053.6665
$ws.Range("J122").Value = 314.5
$ws.Range("K122").Value = 15160.9995
$ws.Range("L122").Value = 943.5
$ws.Range("M122").Value = -12710.9995
$ws.Range("N122").Value = -5843.5
$ws.Range("H129").Value = 715.3333
$ws.Range("I129").Value = 426
$ws.Range("J129").Value = 899.4545000000001
$ws.Range("K129").Value = 1278
$ws.Range("L129").Value = 2698.3635
$ws.Range("M129").Value = 3722
$ws.Range("N129").Value = -12698.3635
$ws.Range("H132").Value = 11121908
$ws.Range("I132").Value = 13340466
$ws.Range("J132").Value = 29119.2
$ws.Range("K132").Value = 40021398
$ws.Range("L132").Value = 87357.60000000001
$ws.Range("M132").Value = -40018868
$ws.Range("N132").Value = -92417.60000000001
$ws.Range("H137").Value = 1870.0785
$ws.Range("I137").Value = 1339.9395
$ws.Range("J137").Value = 2842
$ws.Range("K137").Value = 4019.8185
$ws.Range("L137").Value = 8526
$ws.Range("M137").Value = -1469.8185
$ws.Range("N137").Value = -13626
$ws.Range("H138").Value = 1954.1
$ws.Range("I138").Value = 738.5
$ws.Range("J138").Value = 2185.6428
$ws.Range("K138").Value = 2215.5
$ws.Range("L138").Value = 6556.928400000001
$ws.Range("M138").Value = 2924.5
$ws.Range("N138").Value = -16836.9284
$ws.Range("H141").Value = 1003
$ws.Range("I141").Value = 988.4286
$ws.Range("J141").Value = 1105
$ws.Range("K141").Value = 2965.2858
$ws.Range("L141").Value = 3315
$ws.Range("M141").Value = 2214.7142
$ws.Range("N141").Value = -13675

# Sheet 2 (ARM)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1239
$ws.Range("I2").Value = 810.44446
$ws.Range("K2").Value = 810.44446
$ws.Range("M2").Value = -697.44446
$ws.Range("H116").Value = 1239
$ws.Range("I116").Value = 810.44446
$ws.Range("K116").Value = 810.44446
$ws.Range("M116").Value = 1483.55554
$ws.Range("H122").Value = 2487.0344
$ws.Range("I122").Value = 2441.8076
$ws.Range("J122").Value = 2879
$ws.Range("K122").Value = 7325.4228
$ws.Range("L122").Value = 8637
$ws.Range("M122").Value = -4875.4228
$ws.Range("N122").Value = -13537
$ws.Range("H132").Value = 2654.1628
$ws.Range("I132").Value = 1694.8889
$ws.Range("J132").Value = 4272.9375
$ws.Range("K132").Value = 5084.6667
$ws.Range("L132").Value = 12818.8125
$ws.Range("M132").Value = -2554.6667
$ws.Range("N132").Value = -17878.8125

# Sheet 3 (BSM)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1239
$ws.Range("I3").Value = 810.44446
$ws.Range("K3").Value = 810.44446
$ws.Range("M3").Value = -696.44446
$ws.Range("H59").Value = 53193.332
$ws.Range("I59").Value = 50000
$ws.Range("J59").Value = 54790
$ws.Range("K59").Value = 50000
$ws.Range("L59").Value = 54790
$ws.Range("M59").Value = -49153
$ws.Range("N59").Value = -56484

# Sheet 4 (CRP)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 14287659
$ws.Range("I62").Value = 2003.4445
$ws.Range("J62").Value = 40001840
$ws.Range("K62").Value = 2003.4445
$ws.Range("L62").Value = 40001840
$ws.Range("M62").Value = -1379.4445
$ws.Range("N62").Value = -40003088
$ws.Range("H65").Value = 14287659
$ws.Range("I65").Value = 2003.4445
$ws.Range("J65").Value = 40001840
$ws.Range("K65").Value = 10017.2225
$ws.Range("L65").Value = 200009200
$ws.Range("M65").Value = -6897.2225
$ws.Range("N65").Value = -200015440

# Sheet 5 (CUL)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 486.42426
$ws.Range("I5").Value = 306.96155
$ws.Range("J5").Value = 1153
$ws.Range("K5").Value = 920.88465
$ws.Range("L5").Value = 3459
$ws.Range("M5").Value = -808.88465
$ws.Range("N5").Value = -3683
$ws.Range("H39").Value = 4300.4443
$ws.Range("J39").Value = 4443.4287
$ws.Range("L39").Value = 13330.2861
$ws.Range("N39").Value = -13918.2861
$ws.Range("H55").Value = 2261.5833
$ws.Range("J55").Value = 2945
$ws.Range("L55").Value = 8835
$ws.Range("N55").Value = -9189
$ws.Range("H64").Value = 8115.05
$ws.Range("I64").Value = 20942.2
$ws.Range("J64").Value = 3839.3333
$ws.Range("K64").Value = 62826.60000000001
$ws.Range("L64").Value = 11517.9999
$ws.Range("M64").Value = -62556.60000000001
$ws.Range("N64").Value = -12057.9999
$ws.Range("H67").Value = 8115.05
$ws.Range("I67").Value = 20942.2
$ws.Range("J67").Value = 3839.3333
$ws.Range("K67").Value = 62826.60000000001
$ws.Range("L67").Value = 11517.9999
$ws.Range("M67").Value = -61890.60000000001
$ws.Range("N67").Value = -13389.9999
$ws.Range("H97").Value = 1217.1428
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 1253.3334
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 3760.0002
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -4752.0002
$ws.Range("H131").Value = 14948596
$ws.Range("J131").Value = 29224
$ws.Range("L131").Value = 87672
$ws.Range("N131").Value = -97752
$ws.Range("H135").Value = 486.42426
$ws.Range("I135").Value = 306.96155
$ws.Range("J135").Value = 1153
$ws.Range("K135").Value = 2762.65395
$ws.Range("L135").Value = 10377
$ws.Range("M135").Value = -227.6539499999999
$ws.Range("N135").Value = -15447
$ws.Range("H136").Value = 4949.75
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4949.75
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 14849.25
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -25049.25

# Sheet 6 (GSM)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 3858.9167
$ws.Range("I80").Value = 2240
$ws.Range("K80").Value = 2240
$ws.Range("M80").Value = -1242
$ws.Range("H83").Value = 3858.9167
$ws.Range("I83").Value = 2240
$ws.Range("K83").Value = 11200
$ws.Range("M83").Value = -6208

# Sheet 7 (LTW)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H82").Value = 1969.875
$ws.Range("I82").Value = 1909.3158
$ws.Range("K82").Value = 1909.3158
$ws.Range("M82").Value = -1548.3158
$ws.Range("H85").Value = 1969.875
$ws.Range("I85").Value = 1909.3158
$ws.Range("K85").Value = 1909.3158
$ws.Range("M85").Value = -661.3158000000001
$ws.Range("H88").Value = 19900
$ws.Range("J88").Value = 19900
$ws.Range("L88").Value = 19900
$ws.Range("N88").Value = -20756
$ws.Range("H91").Value = 19900
$ws.Range("J91").Value = 19900
$ws.Range("L91").Value = 19900
$ws.Range("N91").Value = -22864

# Sheet 8 (WVR)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 83335660
$ws.Range("I62").Value = 125002250
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 125002250
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -125001626
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 83335660
$ws.Range("I65").Value = 125002250
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 625011250
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -625008130
$ws.Range("N65").Value = -18740
$ws.Range("H122").Value = 8335242.5
$ws.Range("I122").Value = 9617280
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 28851840
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -28849390
$ws.Range("N122").Value = -10885
